$d = $word.ActiveDocument

function Do-Split($search, $replace) {
    $ok = $d.Content.Find.Execute(
        $search,
        $false, $false, $false, $false, $false, $true, 1, $false,
        $replace,
        2)
    if (-not $ok) {
        Write-Host "WARNING: not found -> $search"
    }
    return $ok
}

# --- Change 1: Portuguese "Programa resumido" paragraph -----------------
# Split "...poluentes atmosféricosA disciplina..." into two <w:t> runs
# separated by a manual line break.
Do-Split `
    "poluentes atmosféricosA disciplina" `
    "poluentes atmosféricos^lA disciplina" | Out-Null

# --- Change 2: English "Programa resumido" paragraph (italic) -----------
# Split "...atmospheric pollutants The discipline..." into two <w:t> runs
# separated by a manual line break (trailing space stays with the first
# part, as in the target diff).
Do-Split `
    "atmospheric pollutants The discipline" `
    "atmospheric pollutants ^lThe discipline" | Out-Null

# --- Change 3: Bibliography paragraph ------------------------------------
# Break the single run of concatenated bibliography entries into many
# <w:t> segments separated by manual line breaks (consecutive entries get
# a double break to form a blank separator line between references).
# Each Find/Replace below targets a short, unique boundary substring
# spanning the end of one bibliography entry and the start of the next,
# so ordering of the calls does not matter (every boundary is unique in
# the original text).

Do-Split `
    "Bibliografia básica:Baird, C.; Cann, M. Quími" `
    "Bibliografia básica:^lBaird, C.; Cann, M. Quími" | Out-Null

Do-Split `
    "okman, 4.ed., 2011. 844p.GUNTER, F.; Introdução ao" `
    "okman, 4.ed., 2011. 844p.^l^lGUNTER, F.; Introdução ao" | Out-Null

Do-Split `
    "Paulo: Editora EPU, 2008.LENZI, E. F.; FAVERO, L.O" `
    "Paulo: Editora EPU, 2008.^l^lLENZI, E. F.; FAVERO, L.O" | Out-Null

Do-Split `
    " Editora LCT, 465p. 2009.Rocha, Julio Cesar; Rosa," `
    " Editora LCT, 465p. 2009.^lRocha, Julio Cesar; Rosa," | Out-Null

Do-Split `
    "legre: Bookman, 2009.  03Seinfeld, J.H. e Pandis, " `
    "legre: Bookman, 2009.  03^l^lSeinfeld, J.H. e Pandis, " | Out-Null

Do-Split `
    " Wiley & Sons Inc., 2006.SPIRO, T. G.; STIGLIANI, " `
    " Wiley & Sons Inc., 2006.^lSPIRO, T. G.; STIGLIANI, " | Out-Null

Do-Split `
    "rentice Hall. 2008. 352p.Bibliografia complementar" `
    "rentice Hall. 2008. 352p.^l^lBibliografia complementar" | Out-Null

Do-Split `
    "ibliografia complementar:Daniel J. Jacob. Introduc" `
    "ibliografia complementar:^lDaniel J. Jacob. Introduc" | Out-Null

Do-Split `
    "n University Press, 1999.FINLAYSON-PITTS, B.J; PIT" `
    "n University Press, 1999.^l^lFINLAYSON-PITTS, B.J; PIT" | Out-Null

Do-Split `
    "go: Academic Press, 2000.Guy P. Brasseur, Max-Plan" `
    "go: Academic Press, 2000.^l^lGuy P. Brasseur, Max-Plan" | Out-Null

Do-Split `
    "mistry 1st Edition, 2017.Jacobson, M.Z. Atmospheri" `
    "mistry 1st Edition, 2017.^lJacobson, M.Z. Atmospheri" | Out-Null

Do-Split `
    "e University Press, 2006.J.H.;MANAHAN, S.E. Enviro" `
    "e University Press, 2006.^lJ.H.;MANAHAN, S.E. Enviro" | Out-Null

Do-Split `
    "L: CRC Press. 753p. 2010.Seinfeld, J.H. e Pandis, " `
    "L: CRC Press. 753p. 2010.^lSeinfeld, J.H. e Pandis, " | Out-Null

Do-Split `
    " Wiley & Sons Inc., 2006.SCHNELLE JR, Karl B; BROW" `
    " Wiley & Sons Inc., 2006.^l^lSCHNELLE JR, Karl B; BROW" | Out-Null

Write-Host "Done"
